$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.512427745349363
$ws.Range("C2").Value = -1.363767227788882
$ws.Range("D2").Value = 59.70675209903689
$ws.Range("E2").Value = -0.0003999063564507079
$ws.Range("F2").Value = 0.6898044308789899
$ws.Range("G2").Value = -0.5004414916667799
$ws.Range("H2").Value = -0.6350004731026235
$ws.Range("I2").Value = 1.621664813878342
$ws.Range("J2").Value = 4.355285082868471
$ws.Range("K2").Value = 33
$ws.Range("L2").Value = -8.080333650738645
$ws.Range("N2").Value = 4.355285082866482
$ws.Range("O2").Value = 5.021951749533149
$ws.Range("B3").Value = 6.780578827585828
$ws.Range("C3").Value = 53.10279935400346
$ws.Range("D3").Value = -89.19799699470447
$ws.Range("E3").Value = -0.00005418064916170746
$ws.Range("F3").Value = 0.37163061194269
$ws.Range("G3").Value = -0.3368040805509991
$ws.Range("H3").Value = 0.0386712405137799
$ws.Range("I3").Value = 1.890841099281444
$ws.Range("J3").Value = 4.355295902663592
$ws.Range("K3").Value = 76
$ws.Range("L3").Value = 64.70060460020855
$ws.Range("N3").Value = 4.355295902660338
$ws.Range("O3").Value = 5.021962569327004
$ws.Range("B4").Value = -0.01513397410389549
$ws.Range("C4").Value = 25.15261976488516
$ws.Range("D4").Value = 44.35061839695146
$ws.Range("E4").Value = 0.2877657036915104
$ws.Range("F4").Value = 1.18443800501441
$ws.Range("G4").Value = -0.5605967116354871
$ws.Range("H4").Value = -0.9456454821981288
$ws.Range("I4").Value = 0.8338599972250136
$ws.Range("J4").Value = 4.355369250924604
$ws.Range("K4").Value = 36
$ws.Range("L4").Value = -5.761773861116499
$ws.Range("N4").Value = 4.355369250931052
$ws.Range("O4").Value = 5.022035917597719
$ws.Range("B5").Value = 5.379649467481933
$ws.Range("C5").Value = -2.599789483648555
$ws.Range("D5").Value = 20.55006959486832
$ws.Range("E5").Value = 50.21933496961683
$ws.Range("F5").Value = 0.5961330041681703
$ws.Range("G5").Value = 0.6765836970136485
$ws.Range("H5").Value = -0.3212336317668543
$ws.Range("I5").Value = -0.3538195067581045
$ws.Range("J5").Value = 4.354102017697073
$ws.Range("K5").Value = 65
$ws.Range("L5").Value = -33.0665014051261
$ws.Range("N5").Value = 4.355385754896673
$ws.Range("O5").Value = 5.02205242156334
$ws.Range("B6").Value = -18.14343379251174
$ws.Range("C6").Value = 49.29212575291513
$ws.Range("D6").Value = 1.855967475239631
$ws.Range("E6").Value = 31.11628356901486
$ws.Range("F6").Value = -0.6542334917483561
$ws.Range("G6").Value = -0.6422895521552101
$ws.Range("H6").Value = 0.489667206518646
$ws.Range("I6").Value = -0.3677553154227562
$ws.Range("J6").Value = 4.354322482717333
$ws.Range("K6").Value = 53
$ws.Range("L6").Value = -18.49407682031511
$ws.Range("N6").Value = 4.355388030443549
$ws.Range("O6").Value = 5.022054697110216
